$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the registration test data (username + email) ---
# A1 held the username "Ashwini12" -> "sonu"
$ws.Range("A1").Value = "sonu"
# C1 held the e-mail "ashwini1234@gmail.com" -> "sonu445@gmail.com"
$ws.Range("C1").Value = "sonu445@gmail.com"

# --- Refresh the mailto: hyperlink on C1 so it points at the new address ---
$ws.Range("C1").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("C1"), "mailto:sonu445@gmail.com")
# Re-apply the built-in Hyperlink cell style (Hyperlinks.Add nudges the style)
$ws.Range("C1").Style = "Hyperlink"

# --- Move the active selection to G7 ---
$ws.Range("G7").Select() | Out-Null
